$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.75
$ws.Range("H2").Value = 3.7
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.38
$ws.Range("L2").Value = 4.75
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.95
$ws.Range("S2").Value = 3.4
$ws.Range("T2").Value = 1.33
$ws.Range("W2").Value = 1.8
$ws.Range("X2").Value = 1.95
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 9.5
$ws.Range("AD2").Value = 29
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 7.5
$ws.Range("AI2").Value = 301
$ws.Range("AK2").Value = 26

$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 2.4
$ws.Range("J3").Value = 3.5
$ws.Range("K3").Value = 2.1
$ws.Range("L3").Value = 3.1
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 2.04
$ws.Range("R3").Value = 1.86
$ws.Range("S3").Value = 3.5
$ws.Range("T3").Value = 1.3
$ws.Range("U3").Value = 1.4
$ws.Range("V3").Value = 2.75
$ws.Range("W3").Value = 1.75
$ws.Range("X3").Value = 2
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 15
$ws.Range("AA3").Value = 11
$ws.Range("AB3").Value = 29
$ws.Range("AC3").Value = 23
$ws.Range("AE3").Value = 10
$ws.Range("AG3").Value = 13
$ws.Range("AH3").Value = 41
$ws.Range("AI3").Value = 201
$ws.Range("AK3").Value = 12
$ws.Range("AL3").Value = 9.5
$ws.Range("AM3").Value = 23
$ws.Range("AN3").Value = 19
$ws.Range("AO3").Value = 29

$ws.Range("G4").Value = 1.29
$ws.Range("H4").Value = 4.75
$ws.Range("J4").Value = 1.8
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.93
$ws.Range("U4").Value = 1.4
$ws.Range("V4").Value = 2.75
$ws.Range("Y4").Value = 5.5
$ws.Range("AE4").Value = 9
$ws.Range("AF4").Value = 9.5
$ws.Range("AR4").Value = 2.5
$ws.Range("AS4").Value = 1.51

$ws.Range("G16").Value = 2.1
$ws.Range("H16").Value = 2.9
$ws.Range("I16").Value = 3.4
$ws.Range("J16").Value = 3
$ws.Range("L16").Value = 4
$ws.Range("W16").Value = 1.91
$ws.Range("X16").Value = 1.8
$ws.Range("Z16").Value = 9.5
$ws.Range("AA16").Value = 9.5
$ws.Range("AB16").Value = 21
$ws.Range("AG16").Value = 15
$ws.Range("AI16").Value = 351
$ws.Range("AJ16").Value = 9
$ws.Range("AK16").Value = 17
$ws.Range("AL16").Value = 13

$ws.Range("G17").Value = 2.3
$ws.Range("I17").Value = 3.1
$ws.Range("J17").Value = 3.2
$ws.Range("M17").Value = 1.08
$ws.Range("N17").Value = 8
$ws.Range("AB17").Value = 23
$ws.Range("AJ17").Value = 8
$ws.Range("AL17").Value = 12

$ws.Range("G18").Value = 1.4
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 1.95
$ws.Range("K18").Value = 2.3
$ws.Range("Q18").Value = 1.9
$ws.Range("R18").Value = 1.95
$ws.Range("AB18").Value = 9
$ws.Range("AJ18").Value = 15
$ws.Range("AK18").Value = 34
$ws.Range("AL18").Value = 21
$ws.Range("AP18").Value = 1.42
$ws.Range("AR18").Value = 2.49
$ws.Range("AS18").Value = 1.54
